$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, pushing existing rows 66-79 down to 67-80.
$ws.Rows("66").Insert()

# Populate the newly inserted row 66 with a new Berenjena price-record.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across this whole data block,
# so copy them straight from the (now shifted) row 67 and just set the
# record-specific values (D,J,K,L,M,O,P) from the diff.
$ws.Range("A66").Value2 = 11
$ws.Range("B66").Value2 = "Vega Monumental Concepción"
$ws.Range("C66").Value2 = "Bíobío"
$ws.Range("D66").Value2 = 44642
$ws.Range("E66").Value2 = 8
$ws.Range("F66").Value2 = 100112001
$ws.Range("G66").Value2 = "Berenjena"
$ws.Range("H66").Value2 = "Sin especificar"
$ws.Range("I66").Value2 = "Primera"
$ws.Range("J66").Value2 = 220
$ws.Range("K66").Value2 = 7000
$ws.Range("L66").Value2 = 7500
$ws.Range("M66").Value2 = 7273
$ws.Range("N66").Value2 = "$/caja 60 unidades"
$ws.Range("O66").Value2 = "Provincia de Chacabuco"
$ws.Range("P66").Value2 = 121
$ws.Range("Q66").Value2 = 60
$ws.Range("R66").Value2 = "Hortaliza"
